$d = $word.ActiveDocument

# The edit removes three whole paragraphs that used to sit right after the
# "LOM3213: Fenomenos de Transporte B (Requisito)" requirement line:
#   1) a blank paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
#
# Use Find to locate the anchor text (robust to any encoding quirks with the
# accented / ©  characters), map the hit back to its paragraph index, then
# delete the three following paragraphs as one contiguous Range so the
# document collapses back together cleanly.

$anchor = $d.Content.Duplicate
$anchor.Find.Execute("LOM3213", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($anchor.Start -ge $p.Range.Start -and $anchor.Start -lt $p.Range.End) {
        $anchorIndex = $i
        break
    }
}

$startPara = $d.Paragraphs.Item($anchorIndex + 1)
$endPara = $d.Paragraphs.Item($anchorIndex + 3)

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
